# users-import-template.xlsx
# Commit: "add column gender, fix bug cc update leave_request"
#
# This sheet is the "Users" import template. The relevant part of the
# change for this workbook is adding a new "gender" column (with a sample
# "male" value) to the example row, and updating the sample employeeId
# value in A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header "gender" in column H (dimension grows from A1:G3 to A1:H3).
$ws.Range("H1").Value = "gender"

# Sample employeeId value changed (EMP260001 -> MAI122).
$ws.Range("A2").Value = "MAI122"

# Sample value for the new "gender" column.
$ws.Range("H2").Value = "male"

# Leave the selection where it ended up after these edits.
$ws.Range("F12").Select()
